$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2024-06-29 Saturday" "2024-06-30 Sunday"

Replace-Text "911÷7=130, 1" "725÷5=145, 0"
Replace-Text "438÷3=146, 0" "510÷5=102, 0"
Replace-Text "605÷9=67, 2" "610÷4=152, 2"
Replace-Text "858÷2=429, 0" "498÷9=55, 3"
Replace-Text "750÷7=107, 1" "261÷2=130, 1"

Replace-Text "878÷6=146, 2" "198÷4=49, 2"
Replace-Text "193÷4=48, 1" "165÷2=82, 1"
Replace-Text "189÷8=23, 5" "979÷6=163, 1"
Replace-Text "503÷8=62, 7" "230÷6=38, 2"
Replace-Text "862÷9=95, 7" "703÷4=175, 3"

Replace-Text "665÷6=110, 5" "882÷6=147, 0"
Replace-Text "686÷9=76, 2" "997÷3=332, 1"
Replace-Text "996÷9=110, 6" "146÷4=36, 2"
Replace-Text "654÷7=93, 3" "137÷7=19, 4"
Replace-Text "965÷7=137, 6" "559÷8=69, 7"

Replace-Text "940÷3=313, 1" "643÷6=107, 1"
Replace-Text "292÷5=58, 2" "284÷2=142, 0"
Replace-Text "942÷6=157, 0" "986÷9=109, 5"
Replace-Text "758÷8=94, 6" "981÷4=245, 1"
Replace-Text "777÷5=155, 2" "679÷9=75, 4"

Replace-Text "820÷5=164, 0" "317÷5=63, 2"
Replace-Text "136÷4=34, 0" "381÷2=190, 1"
Replace-Text "459÷5=91, 4" "327÷7=46, 5"
Replace-Text "582÷2=291, 0" "525÷5=105, 0"
Replace-Text "979÷9=108, 7" "393÷9=43, 6"
